$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 20 de Marzo de 2020 a las 02:46"
$ws.Cells.Item(9, 2).Value = 14267
$ws.Cells.Item(9, 3).Value = 5008
$ws.Cells.Item(9, 5).Value = 13928
$ws.Cells.Item(9, 7).Value = 68
$ws.Cells.Item(9, 8).Value = 218
$ws.Cells.Item(11, 2).Value = 8652
$ws.Cells.Item(11, 3).Value = 239
$ws.Cells.Item(11, 4).Value = 2233
$ws.Cells.Item(11, 5).Value = 6325
$ws.Cells.Item(11, 7).Value = 10
$ws.Cells.Item(11, 8).Value = 94
$ws.Cells.Item(28, 2).Value = 647
$ws.Cells.Item(28, 3).Value = 118
$ws.Cells.Item(28, 5).Value = 638
$ws.Cells.Item(55, 1).Value = "Mexico"
$ws.Cells.Item(55, 2).Value = 164
$ws.Cells.Item(55, 3).Value = 71
$ws.Cells.Item(55, 5).Value = 159
$ws.Cells.Item(55, 6).Value = 1
$ws.Cells.Item(55, 7).Value = 1
$ws.Cells.Item(55, 8).Value = 1
$ws.Cells.Item(56, 1).Value = "Libano"
$ws.Cells.Item(56, 2).Value = 157
$ws.Cells.Item(56, 3).Value = 24
$ws.Cells.Item(56, 4).Value = 4
$ws.Cells.Item(56, 5).Value = 149
$ws.Cells.Item(56, 6).Value = 3
$ws.Cells.Item(56, 8).Value = 4
$ws.Cells.Item(57, 1).Value = "Sudafrica"
$ws.Cells.Item(57, 2).Value = 150
$ws.Cells.Item(57, 3).Value = 34
$ws.Cells.Item(57, 4).Value = 0
$ws.Cells.Item(57, 5).Value = 150
$ws.Cells.Item(57, 6).Value = 0
$ws.Cells.Item(58, 1).Value = "Kuwait"
$ws.Cells.Item(58, 2).Value = 148
$ws.Cells.Item(58, 3).Value = 6
$ws.Cells.Item(58, 4).Value = 18
$ws.Cells.Item(58, 5).Value = 130
$ws.Cells.Item(58, 6).Value = 5
$ws.Cells.Item(58, 8).Value = 0
$ws.Cells.Item(59, 1).Value = "San Marino"
$ws.Cells.Item(59, 2).Value = 144
$ws.Cells.Item(59, 3).Value = 4
$ws.Cells.Item(59, 4).Value = 4
$ws.Cells.Item(59, 5).Value = 126
$ws.Cells.Item(59, 6).Value = 12
$ws.Cells.Item(59, 8).Value = 14
$ws.Cells.Item(60, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(60, 2).Value = 140
$ws.Cells.Item(60, 3).Value = 27
$ws.Cells.Item(60, 4).Value = 31
$ws.Cells.Item(60, 5).Value = 109
$ws.Cells.Item(60, 6).Value = 2
$ws.Cells.Item(60, 8).Value = 0
$ws.Cells.Item(61, 1).Value = "Panama"
$ws.Cells.Item(61, 2).Value = 137
$ws.Cells.Item(61, 3).Value = 28
$ws.Cells.Item(61, 4).Value = 1
$ws.Cells.Item(61, 5).Value = 135
$ws.Cells.Item(61, 6).Value = 7
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 8).Value = 1
$ws.Cells.Item(62, 1).Value = "Argentina"
$ws.Cells.Item(62, 2).Value = 128
$ws.Cells.Item(62, 3).Value = 31
$ws.Cells.Item(62, 4).Value = 3
$ws.Cells.Item(62, 5).Value = 122
$ws.Cells.Item(62, 6).Value = 0
$ws.Cells.Item(62, 7).Value = 1
$ws.Cells.Item(62, 8).Value = 3
$ws.Cells.Item(63, 1).Value = "Eslovaquia"
$ws.Cells.Item(63, 2).Value = 124
$ws.Cells.Item(63, 3).Value = 19
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 124
$ws.Cells.Item(64, 1).Value = "Armenia"
$ws.Cells.Item(64, 2).Value = 122
$ws.Cells.Item(64, 3).Value = 12
$ws.Cells.Item(64, 4).Value = 1
$ws.Cells.Item(64, 5).Value = 121
$ws.Cells.Item(64, 6).Value = 2
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(70, 1).Value = "Uruguay"
$ws.Cells.Item(70, 2).Value = 94
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 94
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(71, 1).Value = "Argelia"
$ws.Cells.Item(71, 2).Value = 90
$ws.Cells.Item(71, 3).Value = 15
$ws.Cells.Item(71, 4).Value = 32
$ws.Cells.Item(71, 5).Value = 49
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 2
$ws.Cells.Item(71, 8).Value = 9
$ws.Cells.Item(72, 1).Value = "Costa Rica"
$ws.Cells.Item(72, 2).Value = 87
$ws.Cells.Item(72, 3).Value = 18
$ws.Cells.Item(72, 4).Value = 0
$ws.Cells.Item(72, 5).Value = 86
$ws.Cells.Item(72, 6).Value = 2
$ws.Cells.Item(72, 8).Value = 1
$ws.Cells.Item(73, 1).Value = "Letonia"
$ws.Cells.Item(73, 2).Value = 86
$ws.Cells.Item(73, 3).Value = 15
$ws.Cells.Item(73, 4).Value = 1
$ws.Cells.Item(73, 5).Value = 85
$ws.Cells.Item(74, 1).Value = "Vietnam"
$ws.Cells.Item(74, 2).Value = 85
$ws.Cells.Item(74, 3).Value = 9
$ws.Cells.Item(74, 4).Value = 16
$ws.Cells.Item(74, 5).Value = 69
$ws.Cells.Item(111, 1).Value = "Cuba"
$ws.Cells.Item(111, 2).Value = 16
$ws.Cells.Item(111, 3).Value = 6
$ws.Cells.Item(111, 8).Value = 1
$ws.Cells.Item(113, 1).Value = "Bolivia"
$ws.Cells.Item(113, 3).Value = 3
$ws.Cells.Item(113, 4).Value = 0
$ws.Cells.Item(113, 5).Value = 15
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(114, 1).Value = "Jamaica"
$ws.Cells.Item(114, 2).Value = 15
$ws.Cells.Item(114, 3).Value = 0
$ws.Cells.Item(114, 4).Value = 2
$ws.Cells.Item(114, 5).Value = 12
$ws.Cells.Item(114, 8).Value = 1
$ws.Cells.Item(115, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(115, 2).Value = 14
$ws.Cells.Item(115, 3).Value = 7
$ws.Cells.Item(115, 5).Value = 14
$ws.Cells.Item(115, 6).Value = 0
$ws.Cells.Item(118, 1).Value = "Paraguay"
$ws.Cells.Item(118, 3).Value = 2
$ws.Cells.Item(118, 6).Value = 1
$ws.Cells.Item(119, 1).Value = "Camerun"
$ws.Cells.Item(119, 2).Value = 13
$ws.Cells.Item(119, 3).Value = 0
$ws.Cells.Item(119, 5).Value = 13
$ws.Cells.Item(120, 1).Value = "Guam"
$ws.Cells.Item(120, 3).Value = 4
$ws.Cells.Item(121, 1).Value = "Honduras"
$ws.Cells.Item(121, 3).Value = 3
$ws.Cells.Item(121, 4).Value = 0
$ws.Cells.Item(121, 5).Value = 12
$ws.Cells.Item(122, 1).Value = "Nigeria"
$ws.Cells.Item(122, 2).Value = 12
$ws.Cells.Item(122, 3).Value = 4
$ws.Cells.Item(122, 4).Value = 1
$ws.Cells.Item(123, 1).Value = "Ruanda"
$ws.Cells.Item(123, 3).Value = 0
$ws.Cells.Item(124, 1).Value = "Ghana"
$ws.Cells.Item(124, 3).Value = 4
$ws.Cells.Item(124, 5).Value = 11
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(128, 1).Value = "Costa de Marfil"
$ws.Cells.Item(128, 3).Value = 0
$ws.Cells.Item(128, 4).Value = 1
$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(129, 1).Value = "Guatemala"
$ws.Cells.Item(129, 3).Value = 1
$ws.Cells.Item(129, 4).Value = 0
$ws.Cells.Item(129, 8).Value = 1
$ws.Cells.Item(131, 1).Value = "Etiopia"
$ws.Cells.Item(131, 3).Value = 1
$ws.Cells.Item(132, 1).Value = "Mauricio"
$ws.Cells.Item(132, 3).Value = 4
$ws.Cells.Item(134, 1).Value = "Mongolia"
$ws.Cells.Item(134, 3).Value = 0
$ws.Cells.Item(136, 1).Value = "Tanzania"
$ws.Cells.Item(136, 3).Value = 3
$ws.Cells.Item(138, 1).Value = "Seychelles"
$ws.Cells.Item(144, 1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(144, 3).Value = 1
$ws.Cells.Item(145, 1).Value = "Namibia"
$ws.Cells.Item(145, 3).Value = 1
$ws.Cells.Item(146, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(147, 1).Value = "San Bartolome"
$ws.Cells.Item(147, 3).Value = 0
$ws.Cells.Item(148, 1).Value = "Kirguistan"
$ws.Cells.Item(148, 3).Value = 0
$ws.Cells.Item(149, 1).Value = "Gabon"
$ws.Cells.Item(149, 3).Value = 0
$ws.Cells.Item(150, 1).Value = "Congo"
$ws.Cells.Item(150, 3).Value = 2
$ws.Cells.Item(153, 1).Value = "Liberia"
$ws.Cells.Item(154, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(155, 1).Value = "Zambia"
$ws.Cells.Item(156, 1).Value = "Santa Lucia"
$ws.Cells.Item(157, 1).Value = "Bermudas"
$ws.Cells.Item(158, 1).Value = "Benin"
$ws.Cells.Item(159, 1).Value = "Groenlandia"
$ws.Cells.Item(160, 1).Value = "Mauritania"
$ws.Cells.Item(162, 1).Value = "El Salvador"
$ws.Cells.Item(163, 1).Value = "Somalia"
$ws.Cells.Item(163, 3).Value = 0
$ws.Cells.Item(164, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(164, 3).Value = 0
$ws.Cells.Item(165, 1).Value = "Surinam"
$ws.Cells.Item(166, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(167, 1).Value = "Butan"
$ws.Cells.Item(168, 1).Value = "Nicaragua"
$ws.Cells.Item(168, 3).Value = 1
$ws.Cells.Item(169, 1).Value = "Santa Sede"
$ws.Cells.Item(169, 3).Value = 0
$ws.Cells.Item(170, 1).Value = "Isla de Man"
$ws.Cells.Item(170, 3).Value = 1
$ws.Cells.Item(174, 1).Value = "Guinea"
$ws.Cells.Item(175, 1).Value = "Republica del Chad"
$ws.Cells.Item(175, 3).Value = 1
$ws.Cells.Item(176, 1).Value = "Niger"
$ws.Cells.Item(177, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(178, 1).Value = "Antigua y Barbuda"
$ws.Cells.Item(179, 1).Value = "Fiyi"
$ws.Cells.Item(180, 1).Value = "Togo"
$ws.Cells.Item(181, 1).Value = "San Vicente y las Granadinas"
